$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell values (recalculated Impact/ASPM stats) ---
$ws.Range("G6").Value = 1.71402632819381
$ws.Range("H6").Value = 9.64841292041411
$ws.Range("I6").Value = 4.8
$ws.Range("M6").Value = 2.624
$ws.Range("N6").Value = 3.9468
$ws.Range("F18").Value = 0.5755
$ws.Range("G18").Value = 0.616445454545455
$ws.Range("L18").Value = 0.4782
$ws.Range("N18").Value = 1.04114
$ws.Range("F19").Value = 0.5755
$ws.Range("G19").Value = 0.616445454545455
$ws.Range("L19").Value = 0.4782
$ws.Range("N19").Value = 1.04114
$ws.Range("G25").Value = 1.48182070716043
$ws.Range("H25").Value = 9.64841292041411
$ws.Range("I25").Value = 4.6
$ws.Range("N25").Value = 3.8164
$ws.Range("G37").Value = 0.669205
$ws.Range("I37").Value = 1.1794
$ws.Range("G38").Value = 0.669205
$ws.Range("I38").Value = 1.1794
$ws.Range("G45").Value = 1.03910184410257
$ws.Range("G48").Value = 1124.20832186581
$ws.Range("H48").Value = 11218.0644260486
$ws.Range("G49").Value = 1124.20832186581
$ws.Range("H49").Value = 11218.0644260486
$ws.Range("G50").Value = 1124.20832186581
$ws.Range("H50").Value = 11218.0644260486
$ws.Range("G51").Value = 1124.20832186581
$ws.Range("H51").Value = 11218.0644260486
$ws.Range("I58").Value = 1.1794
$ws.Range("I59").Value = 1.1794
$ws.Range("G66").Value = 1.01345847583946
$ws.Range("G69").Value = 1159.27611847598
$ws.Range("H69").Value = 11218.0644260486
$ws.Range("G70").Value = 1159.27611847598
$ws.Range("H70").Value = 11218.0644260486
$ws.Range("G71").Value = 1159.27611847598
$ws.Range("H71").Value = 11218.0644260486
$ws.Range("G72").Value = 1159.27611847598
$ws.Range("H72").Value = 11218.0644260486
$ws.Range("G79").Value = 0.717035
$ws.Range("L79").Value = 0.5277500000000001
$ws.Range("G80").Value = 0.717035
$ws.Range("L80").Value = 0.5277500000000001
$ws.Range("G87").Value = 0.866852384785111
$ws.Range("G90").Value = 1213.12791210671
$ws.Range("H90").Value = 11218.0644260486
$ws.Range("G91").Value = 1213.12791210671
$ws.Range("H91").Value = 11218.0644260486
$ws.Range("G92").Value = 1213.12791210671
$ws.Range("H92").Value = 11218.0644260486
$ws.Range("G93").Value = 1213.12791210671
$ws.Range("H93").Value = 11218.0644260486
$ws.Range("G108").Value = 0.7581933021766
$ws.Range("G111").Value = 1596.47417087795
$ws.Range("H111").Value = 11218.0644260486
$ws.Range("I111").Value = 9997.52889
$ws.Range("G112").Value = 1596.47417087795
$ws.Range("H112").Value = 11218.0644260486
$ws.Range("I112").Value = 9997.52889
$ws.Range("G113").Value = 1596.47417087795
$ws.Range("H113").Value = 11218.0644260486
$ws.Range("I113").Value = 9997.52889
$ws.Range("G114").Value = 1596.47417087795
$ws.Range("H114").Value = 11218.0644260486
$ws.Range("I114").Value = 9997.52889
$ws.Range("G129").Value = 0.711784455459927
$ws.Range("G132").Value = 1771.58668517315
$ws.Range("H132").Value = 11218.0644260486
$ws.Range("I132").Value = 9974.64205
$ws.Range("G133").Value = 1771.58668517315
$ws.Range("H133").Value = 11218.0644260486
$ws.Range("I133").Value = 9974.64205
$ws.Range("G134").Value = 1771.58668517315
$ws.Range("H134").Value = 11218.0644260486
$ws.Range("I134").Value = 9974.64205
$ws.Range("G135").Value = 1771.58668517315
$ws.Range("H135").Value = 11218.0644260486
$ws.Range("I135").Value = 9974.64205

# --- Append new rows 148-167 (site-year 2019-2023 attribute rows) ---
# row 148
$ws.Range("A148").Value = "Manawatu at d/s PNCC STP"
$ws.Range("B148").Value = "ASPM"
$ws.Range("C148").Value = "C"
$ws.Range("D148").Value = "2019 - 2023"
$ws.Range("E148").Value = "Impact"
$ws.Range("F148").Value = 0.328
$ws.Range("G148").Value = 0.287
$ws.Range("H148").Value = 0.38
$ws.Range("I148").Value = 0.38
$ws.Range("J148").Value = ""
$ws.Range("K148").Value = ""
$ws.Range("L148").Value = 0.328
$ws.Range("M148").Value = 0.3688
$ws.Range("N148").Value = 0.38
$ws.Range("O148").Value = 1819264.388
$ws.Range("P148").Value = 5525304.917
$ws.Range("Q148").Value = "Palmerston North City"
$ws.Range("R148").Value = "Manawatū"
$ws.Range("S148").Value = "Lower Manawatu"
$ws.Range("T148").Value = "Mana_11a"
$ws.Range("U148").Value = ""

# row 149
$ws.Range("A149").Value = "Manawatu at d/s PNCC STP"
$ws.Range("B149").Value = "Visual Clarity (Sediment class 3)"
$ws.Range("C149").Value = "D"
$ws.Range("D149").Value = "2019 - 2023"
$ws.Range("E149").Value = "Impact"
$ws.Range("F149").Value = 0.31
$ws.Range("G149").Value = 0.681292137337748
$ws.Range("H149").Value = 3.5
$ws.Range("I149").Value = 2.45
$ws.Range("J149").Value = ""
$ws.Range("K149").Value = ""
$ws.Range("L149").Value = 0.4
$ws.Range("M149").Value = 1.55
$ws.Range("N149").Value = 2.0802
$ws.Range("O149").Value = 1819264.388
$ws.Range("P149").Value = 5525304.917
$ws.Range("Q149").Value = "Palmerston North City"
$ws.Range("R149").Value = "Manawatū"
$ws.Range("S149").Value = "Lower Manawatu"
$ws.Range("T149").Value = "Mana_11a"
$ws.Range("U149").Value = "m"

# row 150
$ws.Range("A150").Value = "Manawatu at d/s PNCC STP"
$ws.Range("B150").Value = "DRP (95th Percentile)"
$ws.Range("C150").Value = "C"
$ws.Range("D150").Value = "2019 - 2023"
$ws.Range("E150").Value = "Impact"
$ws.Range("F150").Value = 0.0225
$ws.Range("G150").Value = 0.0234655172413793
$ws.Range("H150").Value = 0.05
$ws.Range("I150").Value = 0.0412
$ws.Range("J150").Value = ""
$ws.Range("K150").Value = ""
$ws.Range("L150").Value = 0.0205
$ws.Range("M150").Value = 0.03164
$ws.Range("N150").Value = 0.03858
$ws.Range("O150").Value = 1819264.388
$ws.Range("P150").Value = 5525304.917
$ws.Range("Q150").Value = "Palmerston North City"
$ws.Range("R150").Value = "Manawatū"
$ws.Range("S150").Value = "Lower Manawatu"
$ws.Range("T150").Value = "Mana_11a"
$ws.Range("U150").Value = "mg/L"

# row 151
$ws.Range("A151").Value = "Manawatu at d/s PNCC STP"
$ws.Range("B151").Value = "DRP (Median)"
$ws.Range("C151").Value = "D"
$ws.Range("D151").Value = "2019 - 2023"
$ws.Range("E151").Value = "Impact"
$ws.Range("F151").Value = 0.0225
$ws.Range("G151").Value = 0.0234655172413793
$ws.Range("H151").Value = 0.05
$ws.Range("I151").Value = 0.0412
$ws.Range("J151").Value = ""
$ws.Range("K151").Value = ""
$ws.Range("L151").Value = 0.0205
$ws.Range("M151").Value = 0.03164
$ws.Range("N151").Value = 0.03858
$ws.Range("O151").Value = 1819264.388
$ws.Range("P151").Value = 5525304.917
$ws.Range("Q151").Value = "Palmerston North City"
$ws.Range("R151").Value = "Manawatū"
$ws.Range("S151").Value = "Lower Manawatu"
$ws.Range("T151").Value = "Mana_11a"
$ws.Range("U151").Value = "mg/L"

# row 152
$ws.Range("A152").Value = "Manawatu at d/s PNCC STP"
$ws.Range("B152").Value = "E coli (>260)"
$ws.Range("C152").Value = "E"
$ws.Range("D152").Value = "2019 - 2023"
$ws.Range("E152").Value = "Impact"
$ws.Range("F152").Value = 400
$ws.Range("G152").Value = 1410.71959913725
$ws.Range("H152").Value = 10157.7367499604
$ws.Range("I152").Value = 8984.799999999999
$ws.Range("J152").Value = 37.9310344827586
$ws.Range("K152").Value = 65.51724137931031
$ws.Range("L152").Value = 360
$ws.Range("M152").Value = 2128
$ws.Range("N152").Value = 5532
$ws.Range("O152").Value = 1819264.388
$ws.Range("P152").Value = 5525304.917
$ws.Range("Q152").Value = "Palmerston North City"
$ws.Range("R152").Value = "Manawatū"
$ws.Range("S152").Value = "Lower Manawatu"
$ws.Range("T152").Value = "Mana_11a"
$ws.Range("U152").Value = "% exceedances over 260/100 mL"

# row 153
$ws.Range("A153").Value = "Manawatu at d/s PNCC STP"
$ws.Range("B153").Value = "E coli (>540)"
$ws.Range("C153").Value = "E"
$ws.Range("D153").Value = "2019 - 2023"
$ws.Range("E153").Value = "Impact"
$ws.Range("F153").Value = 400
$ws.Range("G153").Value = 1410.71959913725
$ws.Range("H153").Value = 10157.7367499604
$ws.Range("I153").Value = 8984.799999999999
$ws.Range("J153").Value = 37.9310344827586
$ws.Range("K153").Value = 65.51724137931031
$ws.Range("L153").Value = 360
$ws.Range("M153").Value = 2128
$ws.Range("N153").Value = 5532
$ws.Range("O153").Value = 1819264.388
$ws.Range("P153").Value = 5525304.917
$ws.Range("Q153").Value = "Palmerston North City"
$ws.Range("R153").Value = "Manawatū"
$ws.Range("S153").Value = "Lower Manawatu"
$ws.Range("T153").Value = "Mana_11a"
$ws.Range("U153").Value = "% exceedances over 540/100 mL"

# row 154
$ws.Range("A154").Value = "Manawatu at d/s PNCC STP"
$ws.Range("B154").Value = "E coli (Median)"
$ws.Range("C154").Value = "E"
$ws.Range("D154").Value = "2019 - 2023"
$ws.Range("E154").Value = "Impact"
$ws.Range("F154").Value = 400
$ws.Range("G154").Value = 1410.71959913725
$ws.Range("H154").Value = 10157.7367499604
$ws.Range("I154").Value = 8984.799999999999
$ws.Range("J154").Value = 37.9310344827586
$ws.Range("K154").Value = 65.51724137931031
$ws.Range("L154").Value = 360
$ws.Range("M154").Value = 2128
$ws.Range("N154").Value = 5532
$ws.Range("O154").Value = 1819264.388
$ws.Range("P154").Value = 5525304.917
$ws.Range("Q154").Value = "Palmerston North City"
$ws.Range("R154").Value = "Manawatū"
$ws.Range("S154").Value = "Lower Manawatu"
$ws.Range("T154").Value = "Mana_11a"
$ws.Range("U154").Value = "E. coli/100 mL"

# row 155
$ws.Range("A155").Value = "Manawatu at d/s PNCC STP"
$ws.Range("B155").Value = "E coli (95th Percentile)"
$ws.Range("C155").Value = "E"
$ws.Range("D155").Value = "2019 - 2023"
$ws.Range("E155").Value = "Impact"
$ws.Range("F155").Value = 400
$ws.Range("G155").Value = 1410.71959913725
$ws.Range("H155").Value = 10157.7367499604
$ws.Range("I155").Value = 8984.799999999999
$ws.Range("J155").Value = 37.9310344827586
$ws.Range("K155").Value = 65.51724137931031
$ws.Range("L155").Value = 360
$ws.Range("M155").Value = 2128
$ws.Range("N155").Value = 5532
$ws.Range("O155").Value = 1819264.388
$ws.Range("P155").Value = 5525304.917
$ws.Range("Q155").Value = "Palmerston North City"
$ws.Range("R155").Value = "Manawatū"
$ws.Range("S155").Value = "Lower Manawatu"
$ws.Range("T155").Value = "Mana_11a"
$ws.Range("U155").Value = "E. coli/100 mL"

# row 156
$ws.Range("A156").Value = "Manawatu at d/s PNCC STP"
$ws.Range("B156").Value = "MCI"
$ws.Range("C156").Value = "D"
$ws.Range("D156").Value = "2019 - 2023"
$ws.Range("E156").Value = "Impact"
$ws.Range("F156").Value = 83.64
$ws.Range("G156").Value = 85.59399999999999
$ws.Range("H156").Value = 101.33
$ws.Range("I156").Value = 101.33
$ws.Range("J156").Value = ""
$ws.Range("K156").Value = ""
$ws.Range("L156").Value = 83.64
$ws.Range("M156").Value = 96.6645
$ws.Range("N156").Value = 101.33
$ws.Range("O156").Value = 1819264.388
$ws.Range("P156").Value = 5525304.917
$ws.Range("Q156").Value = "Palmerston North City"
$ws.Range("R156").Value = "Manawatū"
$ws.Range("S156").Value = "Lower Manawatu"
$ws.Range("T156").Value = "Mana_11a"
$ws.Range("U156").Value = ""

# row 157
$ws.Range("A157").Value = "Manawatu at d/s PNCC STP"
$ws.Range("B157").Value = "Ammoniacal-N (95th Percentile)"
$ws.Range("C157").Value = "C"
$ws.Range("D157").Value = "2019 - 2023"
$ws.Range("E157").Value = "Impact"
$ws.Range("F157").Value = 0.15942
$ws.Range("G157").Value = 0.188419502977187
$ws.Range("H157").Value = 0.796138200862211
$ws.Range("I157").Value = 0.52938
$ws.Range("J157").Value = ""
$ws.Range("K157").Value = ""
$ws.Range("L157").Value = 0.17873
$ws.Range("M157").Value = 0.31521
$ws.Range("N157").Value = 0.43289
$ws.Range("O157").Value = 1819264.388
$ws.Range("P157").Value = 5525304.917
$ws.Range("Q157").Value = "Palmerston North City"
$ws.Range("R157").Value = "Manawatū"
$ws.Range("S157").Value = "Lower Manawatu"
$ws.Range("T157").Value = "Mana_11a"
$ws.Range("U157").Value = "mg NH4-N/L"

# row 158
$ws.Range("A158").Value = "Manawatu at d/s PNCC STP"
$ws.Range("B158").Value = "Ammoniacal-N (Median)"
$ws.Range("C158").Value = "B"
$ws.Range("D158").Value = "2019 - 2023"
$ws.Range("E158").Value = "Impact"
$ws.Range("F158").Value = 0.15942
$ws.Range("G158").Value = 0.188419502977187
$ws.Range("H158").Value = 0.796138200862211
$ws.Range("I158").Value = 0.52938
$ws.Range("J158").Value = ""
$ws.Range("K158").Value = ""
$ws.Range("L158").Value = 0.17873
$ws.Range("M158").Value = 0.31521
$ws.Range("N158").Value = 0.43289
$ws.Range("O158").Value = 1819264.388
$ws.Range("P158").Value = 5525304.917
$ws.Range("Q158").Value = "Palmerston North City"
$ws.Range("R158").Value = "Manawatū"
$ws.Range("S158").Value = "Lower Manawatu"
$ws.Range("T158").Value = "Mana_11a"
$ws.Range("U158").Value = "mg NH4-N/L"

# row 159
$ws.Range("A159").Value = "Manawatu at d/s PNCC STP"
$ws.Range("B159").Value = "Nitrate-N (95th Percentile)"
$ws.Range("C159").Value = "A"
$ws.Range("D159").Value = "2019 - 2023"
$ws.Range("E159").Value = "Impact"
$ws.Range("F159").Value = 0.364
$ws.Range("G159").Value = 0.432327586206897
$ws.Range("H159").Value = 1.06
$ws.Range("I159").Value = 0.927
$ws.Range("J159").Value = ""
$ws.Range("K159").Value = ""
$ws.Range("L159").Value = 0.2435
$ws.Range("M159").Value = 0.7442800000000001
$ws.Range("N159").Value = 0.89872
$ws.Range("O159").Value = 1819264.388
$ws.Range("P159").Value = 5525304.917
$ws.Range("Q159").Value = "Palmerston North City"
$ws.Range("R159").Value = "Manawatū"
$ws.Range("S159").Value = "Lower Manawatu"
$ws.Range("T159").Value = "Mana_11a"
$ws.Range("U159").Value = "mg NO3-N/L"

# row 160
$ws.Range("A160").Value = "Manawatu at d/s PNCC STP"
$ws.Range("B160").Value = "Nitrate-N (Median)"
$ws.Range("C160").Value = "A"
$ws.Range("D160").Value = "2019 - 2023"
$ws.Range("E160").Value = "Impact"
$ws.Range("F160").Value = 0.364
$ws.Range("G160").Value = 0.432327586206897
$ws.Range("H160").Value = 1.06
$ws.Range("I160").Value = 0.927
$ws.Range("J160").Value = ""
$ws.Range("K160").Value = ""
$ws.Range("L160").Value = 0.2435
$ws.Range("M160").Value = 0.7442800000000001
$ws.Range("N160").Value = 0.89872
$ws.Range("O160").Value = 1819264.388
$ws.Range("P160").Value = 5525304.917
$ws.Range("Q160").Value = "Palmerston North City"
$ws.Range("R160").Value = "Manawatū"
$ws.Range("S160").Value = "Lower Manawatu"
$ws.Range("T160").Value = "Mana_11a"
$ws.Range("U160").Value = "mg NO3-N/L"

# row 161
$ws.Range("A161").Value = "Manawatu at d/s PNCC STP"
$ws.Range("B161").Value = "QMCI"
$ws.Range("C161").Value = "C"
$ws.Range("D161").Value = "2019 - 2023"
$ws.Range("E161").Value = "Impact"
$ws.Range("F161").Value = 4.6
$ws.Range("G161").Value = 4.334
$ws.Range("H161").Value = 6.187
$ws.Range("I161").Value = 6.187
$ws.Range("J161").Value = ""
$ws.Range("K161").Value = ""
$ws.Range("L161").Value = 4.6
$ws.Range("M161").Value = 5.71905
$ws.Range("N161").Value = 6.187
$ws.Range("O161").Value = 1819264.388
$ws.Range("P161").Value = 5525304.917
$ws.Range("Q161").Value = "Palmerston North City"
$ws.Range("R161").Value = "Manawatū"
$ws.Range("S161").Value = "Lower Manawatu"
$ws.Range("T161").Value = "Mana_11a"
$ws.Range("U161").Value = ""

# row 162
$ws.Range("A162").Value = "Manawatu at d/s PNCC STP"
$ws.Range("B162").Value = "Soluble Inorganic Nitrogen (95th Percentile)"
$ws.Range("C162").Value = ""
$ws.Range("D162").Value = "2019 - 2023"
$ws.Range("E162").Value = "Impact"
$ws.Range("F162").Value = 0.6885
$ws.Range("G162").Value = 0.701931034482759
$ws.Range("H162").Value = 1.22
$ws.Range("I162").Value = 1.1086
$ws.Range("J162").Value = ""
$ws.Range("K162").Value = ""
$ws.Range("L162").Value = 0.615
$ws.Range("M162").Value = 0.88832
$ws.Range("N162").Value = 1.0672
$ws.Range("O162").Value = 1819264.388
$ws.Range("P162").Value = 5525304.917
$ws.Range("Q162").Value = "Palmerston North City"
$ws.Range("R162").Value = "Manawatū"
$ws.Range("S162").Value = "Lower Manawatu"
$ws.Range("T162").Value = "Mana_11a"
$ws.Range("U162").Value = "g/m3"

# row 163
$ws.Range("A163").Value = "Manawatu at d/s PNCC STP"
$ws.Range("B163").Value = "Soluble Inorganic Nitrogen (Median)"
$ws.Range("C163").Value = ""
$ws.Range("D163").Value = "2019 - 2023"
$ws.Range("E163").Value = "Impact"
$ws.Range("F163").Value = 0.6885
$ws.Range("G163").Value = 0.701931034482759
$ws.Range("H163").Value = 1.22
$ws.Range("I163").Value = 1.1086
$ws.Range("J163").Value = ""
$ws.Range("K163").Value = ""
$ws.Range("L163").Value = 0.615
$ws.Range("M163").Value = 0.88832
$ws.Range("N163").Value = 1.0672
$ws.Range("O163").Value = 1819264.388
$ws.Range("P163").Value = 5525304.917
$ws.Range("Q163").Value = "Palmerston North City"
$ws.Range("R163").Value = "Manawatū"
$ws.Range("S163").Value = "Lower Manawatu"
$ws.Range("T163").Value = "Mana_11a"
$ws.Range("U163").Value = "g/m3"

# row 164
$ws.Range("A164").Value = "Manawatu at d/s PNCC STP"
$ws.Range("B164").Value = "Total Nitrogen (95th Percentile)"
$ws.Range("C164").Value = ""
$ws.Range("D164").Value = "2019 - 2023"
$ws.Range("E164").Value = "Impact"
$ws.Range("F164").Value = 0.925
$ws.Range("G164").Value = 0.996896551724138
$ws.Range("H164").Value = 1.81
$ws.Range("I164").Value = 1.49
$ws.Range("J164").Value = ""
$ws.Range("K164").Value = ""
$ws.Range("L164").Value = 0.875
$ws.Range("M164").Value = 1.2764
$ws.Range("N164").Value = 1.4072
$ws.Range("O164").Value = 1819264.388
$ws.Range("P164").Value = 5525304.917
$ws.Range("Q164").Value = "Palmerston North City"
$ws.Range("R164").Value = "Manawatū"
$ws.Range("S164").Value = "Lower Manawatu"
$ws.Range("T164").Value = "Mana_11a"
$ws.Range("U164").Value = "g/m3"

# row 165
$ws.Range("A165").Value = "Manawatu at d/s PNCC STP"
$ws.Range("B165").Value = "Total Nitrogen (Median)"
$ws.Range("C165").Value = ""
$ws.Range("D165").Value = "2019 - 2023"
$ws.Range("E165").Value = "Impact"
$ws.Range("F165").Value = 0.925
$ws.Range("G165").Value = 0.996896551724138
$ws.Range("H165").Value = 1.81
$ws.Range("I165").Value = 1.49
$ws.Range("J165").Value = ""
$ws.Range("K165").Value = ""
$ws.Range("L165").Value = 0.875
$ws.Range("M165").Value = 1.2764
$ws.Range("N165").Value = 1.4072
$ws.Range("O165").Value = 1819264.388
$ws.Range("P165").Value = 5525304.917
$ws.Range("Q165").Value = "Palmerston North City"
$ws.Range("R165").Value = "Manawatū"
$ws.Range("S165").Value = "Lower Manawatu"
$ws.Range("T165").Value = "Mana_11a"
$ws.Range("U165").Value = "g/m3"

# row 166
$ws.Range("A166").Value = "Manawatu at d/s PNCC STP"
$ws.Range("B166").Value = "Total Phosphorus (95th Percentile)"
$ws.Range("C166").Value = ""
$ws.Range("D166").Value = "2019 - 2023"
$ws.Range("E166").Value = "Impact"
$ws.Range("F166").Value = 0.064
$ws.Range("G166").Value = 0.116913793103448
$ws.Range("H166").Value = 0.969
$ws.Range("I166").Value = 0.3828
$ws.Range("J166").Value = ""
$ws.Range("K166").Value = ""
$ws.Range("L166").Value = 0.054
$ws.Range("M166").Value = 0.16776
$ws.Range("N166").Value = 0.28334
$ws.Range("O166").Value = 1819264.388
$ws.Range("P166").Value = 5525304.917
$ws.Range("Q166").Value = "Palmerston North City"
$ws.Range("R166").Value = "Manawatū"
$ws.Range("S166").Value = "Lower Manawatu"
$ws.Range("T166").Value = "Mana_11a"
$ws.Range("U166").Value = "g/m3"

# row 167
$ws.Range("A167").Value = "Manawatu at d/s PNCC STP"
$ws.Range("B167").Value = "Total Phosphorus (Median)"
$ws.Range("C167").Value = ""
$ws.Range("D167").Value = "2019 - 2023"
$ws.Range("E167").Value = "Impact"
$ws.Range("F167").Value = 0.064
$ws.Range("G167").Value = 0.116913793103448
$ws.Range("H167").Value = 0.969
$ws.Range("I167").Value = 0.3828
$ws.Range("J167").Value = ""
$ws.Range("K167").Value = ""
$ws.Range("L167").Value = 0.054
$ws.Range("M167").Value = 0.16776
$ws.Range("N167").Value = 0.28334
$ws.Range("O167").Value = 1819264.388
$ws.Range("P167").Value = 5525304.917
$ws.Range("Q167").Value = "Palmerston North City"
$ws.Range("R167").Value = "Manawatū"
$ws.Range("S167").Value = "Lower Manawatu"
$ws.Range("T167").Value = "Mana_11a"
$ws.Range("U167").Value = "g/m3"

Write-Host "edit complete"